# Journal de travail Corentin - final commit updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39 (43209): add a comment in column I about using free time for the
# user manual + doc with annexes
$ws.Range("I39").Value = "Durant temps libre manuel d'utilisation + création doc avec annexes"

# Row 40 (43210): fill in the task (C40) for the last journal entry -
# writing the conclusion section of the documentation
$ws.Range("C40").Value = "documentation partie conclusion"

# Row 34 (43188): fix spelling of "malade" -> "maladie" in the task cell
$ws.Range("C34").Value = "maladie"

# Row 40 (43210): comment (F40) for the last journal entry
$ws.Range("F40").Value = "Objectifs, points positifs-négatifs, difficultés, suite, …"

# Move the active selection to reflect where the author left off editing
$ws.Range("C58").Select()
